$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Links")

# Insert a new row at the very top; this shifts all existing rows (1-12) down
# to (2-13), preserving their content, and leaves row 1 empty for new data.
$ws.Rows.Item(1).Insert()

$lorem = 'Lorem Ipsum - это текст-"рыба", часто используемый в печати и вэб-дизайне. Lorem Ipsum является стандартной "рыбой" для текстов на латинице с начала XVI века. В то время некий безымянный печатник создал большую коллекцию размеров и форм шрифтов, используя Lorem Ipsum для распечатки образцов. Lorem Ipsum не только успешно пережил без заметных изменений пять веков, но и перешагнул в электронный дизайн. Его популяризации в новое время послужили публикация листов Letraset с образцами Lorem Ipsum в 60-х годах и, в более недавнее время, программы электронной вёрстки типа Aldus PageMaker, в шаблонах которых используется Lorem Ipsum.'

$ws.Cells.Item(1, 1).Value = 27
$ws.Cells.Item(1, 2).Value = $lorem
$ws.Cells.Item(1, 3).Value = $lorem
$ws.Cells.Item(1, 4).Value = $lorem
